# Update the "Dataset_info" sheet: expand the Town_05 / Town_06 phase
# breakdown (7 and 10 phases respectively, each worth 100 images) and
# shift the Town_07 / Total rows down to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dataset_info")

# --- 1. Insert rows -------------------------------------------------
# Town_05 block (originally rows 33:36 = Phase 1..4) needs 3 more rows
# (Phase 5..7) -> insert 3 blank rows right after row 36 (before the old
# row 37, which is Town_06's first row).
$ws.Rows("37:39").Insert()

# Town_06 block is now at rows 40:42 (old 37:39, still Phase 1..3) and
# needs 7 more rows (Phase 4..10) -> insert 7 blank rows right after the
# new row 42 (before what is now Town_07's first row, 43).
$ws.Rows("43:49").Insert()

# --- 2. Fix up formatting for the newly-inserted blank rows ---------
# Give the 3 new Town_05 rows (37:39) the same border/alignment pattern
# used by every other "middle + bottom" block (copy format only, values
# untouched).
$ws.Range("A37:C37").Copy()
$ws.Range("A37:C38").PasteSpecial(-4122)
$ws.Range("A39:C39").Copy()
$ws.Range("A39:C39").PasteSpecial(-4122)

# Give the 7 new Town_06 rows (43:49) the same pattern: middle rows
# 43:48 then the bottom row 49.
$ws.Range("A43:C43").Copy()
$ws.Range("A43:C48").PasteSpecial(-4122)
$ws.Range("A49:C49").Copy()
$ws.Range("A49:C49").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 3. Write the new row content ------------------------------------
# Town_05 (rows 33-39): Phase 1-7, all 100 images.
$ws.Range("C33").Value = 100
$ws.Range("C34").Value = 100
$ws.Range("C35").Value = 100
$ws.Range("C36").Value = 100
$ws.Range("B37").Value = "Phase 5"
$ws.Range("C37").Value = 100
$ws.Range("B38").Value = "Phase 6"
$ws.Range("C38").Value = 100
$ws.Range("B39").Value = "Phase 7"
$ws.Range("C39").Value = 100

# Town_06 (rows 40-49): Phase 1-10, all 100 images.
$ws.Range("C40").Value = 100
$ws.Range("C41").Value = 100
$ws.Range("C42").Value = 100
$ws.Range("B43").Value = "Phase 4"
$ws.Range("C43").Value = 100
$ws.Range("B44").Value = "Phase 5"
$ws.Range("C44").Value = 100
$ws.Range("B45").Value = "Phase 6"
$ws.Range("C45").Value = 100
$ws.Range("B46").Value = "Phase 7"
$ws.Range("C46").Value = 100
$ws.Range("B47").Value = "Phase 8"
$ws.Range("C47").Value = 100
$ws.Range("B48").Value = "Phase 9"
$ws.Range("C48").Value = 100
$ws.Range("B49").Value = "Phase 10"
$ws.Range("C49").Value = 100

# --- 4. Merge the A-column town labels for the new block sizes ------
$ws.Range("A33:A39").Merge()
$ws.Range("A40:A49").Merge()

# --- 5. Fix the SUM formula / dimension now lives at row 53 ---------
# (the row-insert already shifted "=SUM(C2:C42)" -> "=SUM(C2:C52)" and
# relocated the Total row automatically)

# --- 6. Column C width + view state ---------------------------------
$ws.Columns.Item(3).ColumnWidth = 11.72
$ws.Range("F30").Select()

$excel.CutCopyMode = 0
